$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 11.778396
$ws.Cells.Item(2, 8).Value = 35.335188
$ws.Cells.Item(2, 9).Value = 0.2246658979638982
$ws.Cells.Item(2, 10).Value = 0.2246658979638982
$ws.Cells.Item(2, 13).Value = 34.293805
$ws.Cells.Item(2, 14).Value = 102.881415
$ws.Cells.Item(2, 15).Value = 0.3000373067112135
$ws.Cells.Item(2, 16).Value = 0.3000373067112135
$ws.Cells.Item(2, 17).Value = 403.92601563678
$ws.Cells.Item(2, 18).Value = 3635.33414073102
$ws.Cells.Item(2, 19).Value = 0.06740815093494432
$ws.Cells.Item(2, 20).Value = 0.06740815093494432
$ws.Cells.Item(3, 7).Value = 11.778396
$ws.Cells.Item(3, 8).Value = 35.335188
$ws.Cells.Item(3, 9).Value = 0.2246658979638982
$ws.Cells.Item(3, 10).Value = 0.2246658979638982
$ws.Cells.Item(3, 15).Value = 0.2416702056223797
$ws.Cells.Item(3, 16).Value = 0.2416702056223798
$ws.Cells.Item(3, 17).Value = 325.349151827728
$ws.Cells.Item(3, 18).Value = 2928.142366449552
$ws.Cells.Item(3, 19).Value = 0.05429505375727186
$ws.Cells.Item(3, 20).Value = 0.05429505375727187
$ws.Cells.Item(4, 7).Value = 11.778396
$ws.Cells.Item(4, 8).Value = 35.335188
$ws.Cells.Item(4, 9).Value = 0.2246658979638982
$ws.Cells.Item(4, 10).Value = 0.2246658979638982
$ws.Cells.Item(4, 13).Value = 32.81168366666667
$ws.Cells.Item(4, 14).Value = 98.435051
$ws.Cells.Item(4, 15).Value = 0.2870701922987834
$ws.Cells.Item(4, 16).Value = 0.2870701922987835
$ws.Cells.Item(4, 17).Value = 386.4690036527321
$ws.Cells.Item(4, 18).Value = 3478.221032874588
$ws.Cells.Item(4, 19).Value = 0.06449488253147512
$ws.Cells.Item(4, 20).Value = 0.06449488253147513
$ws.Cells.Item(5, 7).Value = 11.778396
$ws.Cells.Item(5, 8).Value = 35.335188
$ws.Cells.Item(5, 9).Value = 0.2246658979638982
$ws.Cells.Item(5, 10).Value = 0.2246658979638982
$ws.Cells.Item(5, 13).Value = 6.520685
$ws.Cells.Item(5, 14).Value = 19.562055
$ws.Cells.Item(5, 15).Value = 0.05704962646496092
$ws.Cells.Item(5, 16).Value = 0.05704962646496093
$ws.Cells.Item(5, 17).Value = 76.80321012126001
$ws.Cells.Item(5, 18).Value = 691.2288910913401
$ws.Cells.Item(5, 19).Value = 0.01281710555825542
$ws.Cells.Item(5, 20).Value = 0.01281710555825542
$ws.Cells.Item(6, 7).Value = 11.778396
$ws.Cells.Item(6, 8).Value = 35.335188
$ws.Cells.Item(6, 9).Value = 0.2246658979638982
$ws.Cells.Item(6, 10).Value = 0.2246658979638982
$ws.Cells.Item(6, 13).Value = 13.04976133333333
$ws.Cells.Item(6, 14).Value = 39.14928399999999
$ws.Cells.Item(6, 15).Value = 0.1141726689026624
$ws.Cells.Item(6, 16).Value = 0.1141726689026624
$ws.Cells.Item(6, 17).Value = 153.705256689488
$ws.Cells.Item(6, 18).Value = 1383.347310205392
$ws.Cells.Item(6, 19).Value = 0.02565070518195147
$ws.Cells.Item(6, 20).Value = 0.02565070518195148
$ws.Cells.Item(7, 9).Value = 0.2129243426497412
$ws.Cells.Item(7, 10).Value = 0.2129243426497412
$ws.Cells.Item(7, 13).Value = 34.293805
$ws.Cells.Item(7, 14).Value = 102.881415
$ws.Cells.Item(7, 15).Value = 0.3000373067112135
$ws.Cells.Item(7, 16).Value = 0.3000373067112135
$ws.Cells.Item(7, 17).Value = 382.81591526815
$ws.Cells.Item(7, 18).Value = 3445.343237413349
$ws.Cells.Item(7, 19).Value = 0.06388524630188393
$ws.Cells.Item(7, 20).Value = 0.06388524630188393
$ws.Cells.Item(8, 9).Value = 0.2129243426497412
$ws.Cells.Item(8, 10).Value = 0.2129243426497412
$ws.Cells.Item(8, 15).Value = 0.2416702056223797
$ws.Cells.Item(8, 16).Value = 0.2416702056223798
$ws.Cells.Item(8, 19).Value = 0.051457469670173
$ws.Cells.Item(8, 20).Value = 0.05145746967017301
$ws.Cells.Item(9, 9).Value = 0.2129243426497412
$ws.Cells.Item(9, 10).Value = 0.2129243426497412
$ws.Cells.Item(9, 13).Value = 32.81168366666667
$ws.Cells.Item(9, 14).Value = 98.435051
$ws.Cells.Item(9, 15).Value = 0.2870701922987834
$ws.Cells.Item(9, 16).Value = 0.2870701922987835
$ws.Cells.Item(9, 17).Value = 366.2712467847767
$ws.Cells.Item(9, 18).Value = 3296.44122106299
$ws.Cells.Item(9, 19).Value = 0.06112423198955327
$ws.Cells.Item(9, 20).Value = 0.06112423198955328
$ws.Cells.Item(10, 9).Value = 0.2129243426497412
$ws.Cells.Item(10, 10).Value = 0.2129243426497412
$ws.Cells.Item(10, 13).Value = 6.520685
$ws.Cells.Item(10, 14).Value = 19.562055
$ws.Cells.Item(10, 15).Value = 0.05704962646496092
$ws.Cells.Item(10, 16).Value = 0.05704962646496093
$ws.Cells.Item(10, 17).Value = 72.78929813855
$ws.Cells.Item(10, 18).Value = 655.10368324695
$ws.Cells.Item(10, 19).Value = 0.01214725421346508
$ws.Cells.Item(10, 20).Value = 0.01214725421346509
$ws.Cells.Item(11, 9).Value = 0.2129243426497412
$ws.Cells.Item(11, 10).Value = 0.2129243426497412
$ws.Cells.Item(11, 13).Value = 13.04976133333333
$ws.Cells.Item(11, 14).Value = 39.14928399999999
$ws.Cells.Item(11, 15).Value = 0.1141726689026624
$ws.Cells.Item(11, 16).Value = 0.1141726689026624
$ws.Cells.Item(11, 17).Value = 145.6722673045733
$ws.Cells.Item(11, 18).Value = 1311.05040574116
$ws.Cells.Item(11, 19).Value = 0.02431014047466593
$ws.Cells.Item(11, 20).Value = 0.02431014047466594
$ws.Cells.Item(12, 7).Value = 13.44231133333333
$ws.Cells.Item(12, 8).Value = 40.326934
$ws.Cells.Item(12, 9).Value = 0.256404093257997
$ws.Cells.Item(12, 10).Value = 0.256404093257997
$ws.Cells.Item(12, 13).Value = 34.293805
$ws.Cells.Item(12, 14).Value = 102.881415
$ws.Cells.Item(12, 15).Value = 0.3000373067112135
$ws.Cells.Item(12, 16).Value = 0.3000373067112135
$ws.Cells.Item(12, 17).Value = 460.9880036146234
$ws.Cells.Item(12, 18).Value = 4148.892032531609
$ws.Cells.Item(12, 19).Value = 0.07693079357086023
$ws.Cells.Item(12, 20).Value = 0.07693079357086023
$ws.Cells.Item(13, 7).Value = 13.44231133333333
$ws.Cells.Item(13, 8).Value = 40.326934
$ws.Cells.Item(13, 9).Value = 0.256404093257997
$ws.Cells.Item(13, 10).Value = 0.256404093257997
$ws.Cells.Item(13, 15).Value = 0.2416702056223797
$ws.Cells.Item(13, 16).Value = 0.2416702056223798
$ws.Cells.Item(13, 17).Value = 371.3107108051262
$ws.Cells.Item(13, 18).Value = 3341.796397246136
$ws.Cells.Item(13, 19).Value = 0.06196522994007997
$ws.Cells.Item(13, 20).Value = 0.06196522994007998
$ws.Cells.Item(14, 7).Value = 13.44231133333333
$ws.Cells.Item(14, 8).Value = 40.326934
$ws.Cells.Item(14, 9).Value = 0.256404093257997
$ws.Cells.Item(14, 10).Value = 0.256404093257997
$ws.Cells.Item(14, 13).Value = 32.81168366666667
$ws.Cells.Item(14, 14).Value = 98.435051
$ws.Cells.Item(14, 15).Value = 0.2870701922987834
$ws.Cells.Item(14, 16).Value = 0.2870701922987835
$ws.Cells.Item(14, 17).Value = 441.0648672181816
$ws.Cells.Item(14, 18).Value = 3969.583804963634
$ws.Cells.Item(14, 19).Value = 0.0736059723577684
$ws.Cells.Item(14, 20).Value = 0.07360597235776842
$ws.Cells.Item(15, 7).Value = 13.44231133333333
$ws.Cells.Item(15, 8).Value = 40.326934
$ws.Cells.Item(15, 9).Value = 0.256404093257997
$ws.Cells.Item(15, 10).Value = 0.256404093257997
$ws.Cells.Item(15, 13).Value = 6.520685
$ws.Cells.Item(15, 14).Value = 19.562055
$ws.Cells.Item(15, 15).Value = 0.05704962646496092
$ws.Cells.Item(15, 16).Value = 0.05704962646496093
$ws.Cells.Item(15, 17).Value = 87.65307787659668
$ws.Cells.Item(15, 18).Value = 788.8777008893701
$ws.Cells.Item(15, 19).Value = 0.01462775774445573
$ws.Cells.Item(15, 20).Value = 0.01462775774445574
$ws.Cells.Item(16, 7).Value = 13.44231133333333
$ws.Cells.Item(16, 8).Value = 40.326934
$ws.Cells.Item(16, 9).Value = 0.256404093257997
$ws.Cells.Item(16, 10).Value = 0.256404093257997
$ws.Cells.Item(16, 13).Value = 13.04976133333333
$ws.Cells.Item(16, 14).Value = 39.14928399999999
$ws.Cells.Item(16, 15).Value = 0.1141726689026624
$ws.Cells.Item(16, 16).Value = 0.1141726689026624
$ws.Cells.Item(16, 17).Value = 175.4189546683617
$ws.Cells.Item(16, 18).Value = 1578.770592015256
$ws.Cells.Item(16, 19).Value = 0.02927433964483265
$ws.Cells.Item(16, 20).Value = 0.02927433964483266
$ws.Cells.Item(17, 7).Value = 3.363226
$ws.Cells.Item(17, 8).Value = 10.089678
$ws.Cells.Item(17, 9).Value = 0.06415153551854848
$ws.Cells.Item(17, 10).Value = 0.06415153551854848
$ws.Cells.Item(17, 13).Value = 34.293805
$ws.Cells.Item(17, 14).Value = 102.881415
$ws.Cells.Item(17, 15).Value = 0.3000373067112135
$ws.Cells.Item(17, 16).Value = 0.3000373067112135
$ws.Cells.Item(17, 17).Value = 115.33781661493
$ws.Cells.Item(17, 18).Value = 1038.04034953437
$ws.Cells.Item(17, 19).Value = 0.01924785393837403
$ws.Cells.Item(17, 20).Value = 0.01924785393837403
$ws.Cells.Item(18, 7).Value = 3.363226
$ws.Cells.Item(18, 8).Value = 10.089678
$ws.Cells.Item(18, 9).Value = 0.06415153551854848
$ws.Cells.Item(18, 10).Value = 0.06415153551854848
$ws.Cells.Item(18, 15).Value = 0.2416702056223797
$ws.Cells.Item(18, 16).Value = 0.2416702056223798
$ws.Cells.Item(18, 17).Value = 92.90082677683466
$ws.Cells.Item(18, 18).Value = 836.1074409915119
$ws.Cells.Item(18, 19).Value = 0.01550351477975901
$ws.Cells.Item(18, 20).Value = 0.01550351477975901
$ws.Cells.Item(19, 7).Value = 3.363226
$ws.Cells.Item(19, 8).Value = 10.089678
$ws.Cells.Item(19, 9).Value = 0.06415153551854848
$ws.Cells.Item(19, 10).Value = 0.06415153551854848
$ws.Cells.Item(19, 13).Value = 32.81168366666667
$ws.Cells.Item(19, 14).Value = 98.435051
$ws.Cells.Item(19, 15).Value = 0.2870701922987834
$ws.Cells.Item(19, 16).Value = 0.2870701922987835
$ws.Cells.Item(19, 17).Value = 110.3531076115087
$ws.Cells.Item(19, 18).Value = 993.177968503578
$ws.Cells.Item(19, 19).Value = 0.01841599363757195
$ws.Cells.Item(19, 20).Value = 0.01841599363757195
$ws.Cells.Item(20, 7).Value = 3.363226
$ws.Cells.Item(20, 8).Value = 10.089678
$ws.Cells.Item(20, 9).Value = 0.06415153551854848
$ws.Cells.Item(20, 10).Value = 0.06415153551854848
$ws.Cells.Item(20, 13).Value = 6.520685
$ws.Cells.Item(20, 14).Value = 19.562055
$ws.Cells.Item(20, 15).Value = 0.05704962646496092
$ws.Cells.Item(20, 16).Value = 0.05704962646496093
$ws.Cells.Item(20, 17).Value = 21.93053732981
$ws.Cells.Item(20, 18).Value = 197.37483596829
$ws.Cells.Item(20, 19).Value = 0.003659821138486864
$ws.Cells.Item(20, 20).Value = 0.003659821138486865
$ws.Cells.Item(21, 7).Value = 3.363226
$ws.Cells.Item(21, 8).Value = 10.089678
$ws.Cells.Item(21, 9).Value = 0.06415153551854848
$ws.Cells.Item(21, 10).Value = 0.06415153551854848
$ws.Cells.Item(21, 13).Value = 13.04976133333333
$ws.Cells.Item(21, 14).Value = 39.14928399999999
$ws.Cells.Item(21, 15).Value = 0.1141726689026624
$ws.Cells.Item(21, 16).Value = 0.1141726689026624
$ws.Cells.Item(21, 17).Value = 43.88929661006132
$ws.Cells.Item(21, 18).Value = 395.0036694905519
$ws.Cells.Item(21, 19).Value = 0.007324352024356619
$ws.Cells.Item(21, 20).Value = 0.007324352024356621
$ws.Cells.Item(22, 7).Value = 12.679511
$ws.Cells.Item(22, 8).Value = 38.038533
$ws.Cells.Item(22, 9).Value = 0.2418541306098152
$ws.Cells.Item(22, 10).Value = 0.2418541306098152
$ws.Cells.Item(22, 13).Value = 34.293805
$ws.Cells.Item(22, 14).Value = 102.881415
$ws.Cells.Item(22, 15).Value = 0.3000373067112135
$ws.Cells.Item(22, 16).Value = 0.3000373067112135
$ws.Cells.Item(22, 17).Value = 434.828677729355
$ws.Cells.Item(22, 18).Value = 3913.458099564195
$ws.Cells.Item(22, 19).Value = 0.07256526196515099
$ws.Cells.Item(22, 20).Value = 0.07256526196515099
$ws.Cells.Item(23, 7).Value = 12.679511
$ws.Cells.Item(23, 8).Value = 38.038533
$ws.Cells.Item(23, 9).Value = 0.2418541306098152
$ws.Cells.Item(23, 10).Value = 0.2418541306098152
$ws.Cells.Item(23, 15).Value = 0.2416702056223797
$ws.Cells.Item(23, 16).Value = 0.2416702056223798
$ws.Cells.Item(23, 17).Value = 350.2402321538813
$ws.Cells.Item(23, 18).Value = 3152.162089384932
$ws.Cells.Item(23, 19).Value = 0.05844893747509591
$ws.Cells.Item(23, 20).Value = 0.05844893747509592
$ws.Cells.Item(24, 7).Value = 12.679511
$ws.Cells.Item(24, 8).Value = 38.038533
$ws.Cells.Item(24, 9).Value = 0.2418541306098152
$ws.Cells.Item(24, 10).Value = 0.2418541306098152
$ws.Cells.Item(24, 13).Value = 32.81168366666667
$ws.Cells.Item(24, 14).Value = 98.435051
$ws.Cells.Item(24, 15).Value = 0.2870701922987834
$ws.Cells.Item(24, 16).Value = 0.2870701922987835
$ws.Cells.Item(24, 17).Value = 416.0361039800204
$ws.Cells.Item(24, 18).Value = 3744.324935820183
$ws.Cells.Item(24, 19).Value = 0.06942911178241473
$ws.Cells.Item(24, 20).Value = 0.06942911178241475
$ws.Cells.Item(25, 7).Value = 12.679511
$ws.Cells.Item(25, 8).Value = 38.038533
$ws.Cells.Item(25, 9).Value = 0.2418541306098152
$ws.Cells.Item(25, 10).Value = 0.2418541306098152
$ws.Cells.Item(25, 13).Value = 6.520685
$ws.Cells.Item(25, 14).Value = 19.562055
$ws.Cells.Item(25, 15).Value = 0.05704962646496092
$ws.Cells.Item(25, 16).Value = 0.05704962646496093
$ws.Cells.Item(25, 17).Value = 82.679097185035
$ws.Cells.Item(25, 18).Value = 744.1118746653151
$ws.Cells.Item(25, 19).Value = 0.01379768781029783
$ws.Cells.Item(25, 20).Value = 0.01379768781029783
$ws.Cells.Item(26, 7).Value = 12.679511
$ws.Cells.Item(26, 8).Value = 38.038533
$ws.Cells.Item(26, 9).Value = 0.2418541306098152
$ws.Cells.Item(26, 10).Value = 0.2418541306098152
$ws.Cells.Item(26, 13).Value = 13.04976133333333
$ws.Cells.Item(26, 14).Value = 39.14928399999999
$ws.Cells.Item(26, 15).Value = 0.1141726689026624
$ws.Cells.Item(26, 16).Value = 0.1141726689026624
$ws.Cells.Item(26, 17).Value = 165.4645923733746
$ws.Cells.Item(26, 18).Value = 1489.181331360372
$ws.Cells.Item(26, 19).Value = 0.02761313157685568
$ws.Cells.Item(26, 20).Value = 0.02761313157685569
